$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 113, shifting existing rows 113-174 down to 114-175
$ws.Rows("113:113").Insert()

# Populate the newly inserted row 113 with the new price entry
$ws.Cells.Item(113, 1).Value = 5
$ws.Cells.Item(113, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(113, 3).Value = "Maule"
$ws.Cells.Item(113, 4).Value = 44518
$ws.Cells.Item(113, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(113, 5).Value = 7
$ws.Cells.Item(113, 6).Value = "Fruta"
$ws.Cells.Item(113, 7).Value = 100108
$ws.Cells.Item(113, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(113, 9).Value = 100108005
$ws.Cells.Item(113, 10).Value = "Piña"
$ws.Cells.Item(113, 11).Value = "Caramelo"
$ws.Cells.Item(113, 12).Value = "Tercera"
$ws.Cells.Item(113, 13).Value = 200
$ws.Cells.Item(113, 14).Value = 18000
$ws.Cells.Item(113, 15).Value = 18000
$ws.Cells.Item(113, 16).Value = 18000
$ws.Cells.Item(113, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(113, 18).Value = "Ecuador"
$ws.Cells.Item(113, 19).Value = 1125
$ws.Cells.Item(113, 20).Value = 16
